$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2020" column (N) was added to the table, mirroring the layout/
# formatting of the existing "2019" column (M). Copy M's per-row cell
# formatting into N first so the new column's styles (fonts, number
# formats, borders, alignment) line up with the rest of the table.
$ws.Range("M4:M16").Copy() | Out-Null
$ws.Range("N4:N16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Header
$ws.Range("N4").Value = 2020

# Data rows
$ws.Range("N5").Value = 588.70000000000005
$ws.Range("N6").Value = 62.2
$ws.Range("N7").Value = 99.4
$ws.Range("N8").Value = 6.1
$ws.Range("N9").Value = "-"
$ws.Range("N10").Value = 71
$ws.Range("N11").Value = 136.30000000000001
$ws.Range("N12").Value = 103.3
$ws.Range("N13").Value = 103.2
$ws.Range("N14").Value = 1.8
$ws.Range("N15").Value = "-"
$ws.Range("N16").Value = 5.4

# N10 carries its own one-decimal number format, distinct from the style
# copied from M10.
$ws.Range("N10").NumberFormat = "0.0"

# Match the saved selection state.
$ws.Range("P15").Select() | Out-Null
